# Apply the "edited request command messages" edit to the single slide of
# the HighLevelSequenceDiagrams deck.
#
# Summary of content changes:
#   - "/order delete 1"                -> "delete 1"                  (TextBox 23, also reflowed/resized)
#   - execute("/order delete 1")       -> execute("delete 1")         (TextBox 25, also repositioned/resized)
#   - deleteOrder(o)                   -> deletePerson(p)             (TextBox 28)
#   - post(OrderBookChangedEvent)      -> post(AddressBookChangedEvent)           (TextBox 32 and TextBox 61)
#   - handleOrderBookChangedEvent()    -> handleAddresssBookChangedEvent()        (TextBox 73 and TextBox 49)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Smart/curly quote characters used by "execute(...)" label.
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”

# PowerPoint shape geometry is stored in EMU (1 inch = 914400 EMU) but the
# Shape.Left/Top/Width/Height COM properties are expressed in points (1 inch
# = 72 pt) as single-precision floats. A plain emu/914400*72 conversion can
# therefore land 1 EMU away from the intended integer once it round-trips
# through float32, so search the neighbourhood of the naive conversion for
# the closest point value that truncates back to exactly the target EMU.
function EMUToPoints([double]$emuTarget) {
    $base = $emuTarget / 914400.0 * 72.0
    $bestPt = $base
    $bestDelta = 999999.0
    for ($i = -400; $i -le 400; $i++) {
        $delta = 0.0000005 * $i
        $pt = $base + $delta
        $f32 = [float]$pt
        $emuF32 = [double]$f32 * 914400.0 / 72.0
        $trunc = [int64]$emuF32
        if ($trunc -eq $emuTarget) {
            $adelta = [math]::Abs($delta)
            if ($adelta -lt $bestDelta) {
                $bestDelta = $adelta
                $bestPt = $pt
            }
        }
    }
    return $bestPt
}

# ---------------------------------------------------------------------------
# 1) TextBox 23: "/order delete 1" -> "delete 1", shrink/move the box.
# ---------------------------------------------------------------------------
$tb23 = $s.Shapes.Item("TextBox 23")
$tb23.TextFrame.TextRange.Text = "delete 1"
$tb23.Left  = EMUToPoints 466818
$tb23.Width = EMUToPoints 860170

# ---------------------------------------------------------------------------
# 2) TextBox 25: execute("/order delete 1") -> execute("delete 1"), reposition.
# ---------------------------------------------------------------------------
$tb25 = $s.Shapes.Item("TextBox 25")
$tb25.TextFrame.TextRange.Text = "execute(" + $ldq + "delete 1" + $rdq + ")"
$tb25.Left  = EMUToPoints 2166172
$tb25.Top   = EMUToPoints 1453379
$tb25.Width = EMUToPoints 1424846

# ---------------------------------------------------------------------------
# 3) TextBox 28: deleteOrder(o) -> deletePerson(p) (two runs, keep formatting).
# ---------------------------------------------------------------------------
$tb28 = $s.Shapes.Item("TextBox 28")
$tr28 = $tb28.TextFrame.TextRange
$tr28.Characters(1, 11).Text = "deletePerson"   # "deleteOrder" -> "deletePerson"
$tr28.Characters(13, 3).Text = "(p)"            # "(o)" -> "(p)"

# ---------------------------------------------------------------------------
# 4) TextBox 32 & TextBox 61: post(OrderBookChangedEvent) -> post(AddressBookChangedEvent)
# ---------------------------------------------------------------------------
foreach ($name in @("TextBox 32", "TextBox 61")) {
    $sh = $s.Shapes.Item($name)
    $tr = $sh.TextFrame.TextRange
    $tr.Characters(6, 21).Text = "AddressBookChangedEvent"   # "OrderBookChangedEvent" -> "AddressBookChangedEvent"
}

# ---------------------------------------------------------------------------
# 5) TextBox 73 & TextBox 49: handleOrderBookChangedEvent() -> handleAddresssBookChangedEvent()
#    These boxes auto-fit their height to the (fixed-width) wrapped text; the
#    slightly longer replacement text still fits on a single line in the
#    original deck, so pin the height back to its original value (215444 EMU)
#    after the text edit in case the text-layout re-flow nudges it.
# ---------------------------------------------------------------------------
foreach ($name in @("TextBox 73", "TextBox 49")) {
    $sh = $s.Shapes.Item($name)
    $tr = $sh.TextFrame.TextRange
    $tr.Characters(1, 27).Text = "handleAddresssBookChangedEvent"   # "handleOrderBookChangedEvent" -> "handleAddresssBookChangedEvent"
    $sh.Height = EMUToPoints 215444
}
